$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the "datetimeFigureOut" footer field cached text (10/7/20 ->
#    12/10/20) on the slide master and every slide layout, mirroring what
#    PowerPoint does when it recomputes the date placeholder's cached text.
# ---------------------------------------------------------------------------
$newDate = "12/10/20"
$ppPlaceholderDate = 16

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDate = $true
            }
        } catch {
        }
        if ($isDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $isDate = $false
            try {
                if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDate = $true
                }
            } catch {
            }
            if ($isDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 9 ("5 weeks before the FRM (Jan 18)") - append a new bullet to the
#    "Content Placeholder 2" shape, at the top outline level, right after the
#    existing "Work with Bill ..." sub-bullet.
# ---------------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)

$contentShape = $null
for ($i = 1; $i -le $slide9.Shapes.Count; $i++) {
    $shp = $slide9.Shapes.Item($i)
    if ($shp.Name -eq "Content Placeholder 2") {
        $contentShape = $shp
    }
}

$tr = $contentShape.TextFrame.TextRange

$line1 = "Resolve all remaining open issues from the previous lists"
$line2 = "Work with Bill to merge any changes (prior to errata vote at the FRM)"
$line3 = "All PRs need to be posted by Jan 11 to give time for Bill to review and merge them."

$tr.Text = $line1 + "`r" + $line2 + "`r" + $line3

# Restore the original outline level (2nd level) on the "Work with Bill ..."
# paragraph; the newly appended paragraph stays at the default (top) level.
$start2 = $line1.Length + 2
$para2 = $tr.Characters($start2, $line2.Length)
$para2.IndentLevel = 2
